# Regenerate save_data to use K (strikeouts) instead of Strike# for the
# hill_rich 2022 sheet: recompute column G (K) values and rewrite them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values (column G), rows 2-30, computed from the regenerated
# box-score derived stats (K replaces the old Strike# count).
$kValues = @{
    2  = 6
    3  = 9
    4  = 4
    5  = 4
    6  = 7
    7  = 3
    8  = 2
    9  = 11
    10 = 4
    11 = 0
    12 = 1
    13 = 3
    14 = 5
    15 = 6
    16 = 4
    17 = 6
    18 = 5
    19 = 5
    20 = 1
    21 = 2
    22 = 4
    23 = 6
    24 = 4
    25 = 1
    26 = 2
    27 = 4
    28 = 1
    29 = 3
    30 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
